$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.468.65'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '2.286.44'
$ws.Range("E3").Value = '  +1.26%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.64%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0958'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("E10").Value = '  +0.80%  '

$ws.Range("E11").Value = '  +4.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.73'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").Value = '2.693.82'
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.93%  '

$ws.Range("D15").Value = '54.429.86'
$ws.Range("E15").Value = '  +0.16%  '

$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '2.302.11'
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.82%  '

$ws.Range("E19").Value = '  +3.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '304.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.76%  '

$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("E25").Value = '  +1.36%  '

$ws.Range("E26").Value = '  +3.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("E28").Value = '  +1.78%  '

$ws.Range("D29").Value = '0.0₃0694'
$ws.Range("E29").Value = '  +1.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.79%  '

$ws.Range("E31").Value = '  +1.48%  '

$ws.Range("E32").Value = '  -0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.966'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("E37").Value = '  +2.84%  '

$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.42'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.07'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0496'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.87%  '

$ws.Range("E44").Value = '  +0.80%  '

$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '242.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.28%  '

$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("E48").Value = '  +1.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.43'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.08%  '

$ws.Range("E51").Value = '  -0.44%  '
